# "Adding test cases for watch list(TestCase_E5,TestCase_E6)"
#
# The watch-list test cases (TestCase_E5/E6, rows 6-7) already existed in the
# sheet but were marked as not-yet-run; this commit flips the three earlier
# cases (TestCase_E1..E3, rows 2-4) from SKIP to PASS in the Results column
# now that they've actually been executed. Excel drops the "SKIP" shared
# string automatically once nothing references it any more.
#
# Note: the diff also shows the bookViews/workbookView windowWidth and
# windowHeight shrinking (15150x10125 -> 14310x5130). That is pure window-
# chrome state carried over verbatim from however the workbook was last
# saved on the author's machine - it is not backed by any Workbook/Window
# COM property in this host (Application/ActiveWindow Width & Height are
# writable in-session but are not part of the persisted workbook view and
# do not round-trip into bookViews on save), so there is no COM call here
# that can reproduce it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Results column (E): SKIP -> PASS for TestCase_E1, TestCase_E2, TestCase_E3.
$ws.Range("E2").Value = "PASS"
$ws.Range("E3").Value = "PASS"
$ws.Range("E4").Value = "PASS"

# The active selection on the sheet moves from C6 to C7.
$ws.Range("C7").Select()
